{"js": "// UI UPDATE - ADD RAPPEL TO BUDGET RESUME\n// Update the three \"Par Cr\u00e9ance / Par mandats / Par Chapitre\" amounts in\n// the budget resume table for each of the two changed rows.\nconst replacements = [\n  [\"41 940 000,00\", \"41 560 000,00\"],\n  [\"3 774 600,00\", \"3 740 400,00\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# UI UPDATE - ADD RAPPEL TO BUDGET RESUME\n# Update the three \"Par Cr\u00e9ance / Par mandats / Par Chapitre\" amounts in\n# the budget resume table for each of the two changed rows.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$pairs = @(\n    @{ Old = \"41 940 000,00\"; New = \"41 560 000,00\" },\n    @{ Old = \"3 774 600,00\";  New = \"3 740 400,00\" }\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        try {\n            $cell = $t.Cell($r, $c)\n        } catch {\n            continue\n        }\n        $cellRange = $cell.Range\n        $cellText = $cellRange.Text\n        foreach ($pair in $pairs) {\n            if ($cellText -like \"*$($pair.Old)*\") {\n                $cellRange.Text = $pair.New\n            }\n        }\n    }\n}\n"}
